$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: fill in the previously-empty cedula (B3) for Jose Alberto Molina
$ws.Range("B3").Value = "161-090392-0001F"

# Row 4: new student record - Cristina Maria Leiva Fajardo
$ws.Range("A4").Value = "2009-30746"
$ws.Range("B4").Value = "171-090392-0001Q"
$ws.Range("C4").Value = "Cristina Maria"
$ws.Range("D4").Value = "Leiva Fajardo"
$ws.Range("E4").Value = "F"
# Use the raw date serial number so no time-of-day fraction is attached;
# the cell already carries the d-m-yyyy custom number format.
$ws.Range("F4").Value = 33765
$ws.Range("G4").Value = "leivafajardo@gmail.com"
$ws.Range("J4").Value = "La concha de la lora"
$ws.Range("K4").Value = "Ingeniería de sistemas"

# Hyperlink the new e-mail address, matching the existing G3 hyperlink
$ws.Hyperlinks.Add($ws.Range("G4"), "mailto:leivafajardo@gmail.com")

# Adding the hyperlink re-styles G4 with a generic hyperlink style; restore
# the original left-aligned hyperlink style already used by G3.
$ws.Range("G3").Copy()
$ws.Range("G4").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Move the active selection, as recorded in the saved view state
$ws.Range("B9").Select() | Out-Null
